# =====================================================================
# Edit: add 2022-Q1 data
#
# Before:  sheets = [2020-Q4, 2021-Q1, 2021-Q2, 2021-Q3, 2021-Q4, 总计]
# After:   sheets = [2020-Q4, 2021-Q1, 2021-Q2, 2021-Q3, 2021-Q4, 2022-Q1, 总计]
#
# The sheet that used to be named "总计" (sheetId 6) is repurposed to hold
# the per-fund holding detail for 2022-Q1 (same shape as the other quarter
# sheets). A brand new "总计" sheet (sheetId 7) is appended at the end,
# holding the historical summary table with a new 2022-Q1 row prepended.
# =====================================================================

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Reference sheets by their (stable) index before we start mutating names.
# ---------------------------------------------------------------------
$wsQ4_2021   = $wb.Worksheets.Item(5)   # "2021-Q4" -- used as a style donor
$wsOldTotal  = $wb.Worksheets.Item(6)   # currently named "总计"

# ---------------------------------------------------------------------
# Step 1: the old "总计" sheet becomes the "2022-Q1" fund-detail sheet.
# ---------------------------------------------------------------------
$wsOldTotal.Name = "2022-Q1"
$ws = $wsOldTotal

# Wipe the old 4-column summary content (A1:D6) entirely.
$ws.Range("A1:D6").Clear()

# ---- header row (style copied from the "2021-Q4" sheet's header cell) ----
$ws.Range("B1").Value = '基金代码'
$ws.Range("C1").Value = '基金名称'
$ws.Range("D1").Value = '基金规模'
$ws.Range("E1").Value = '股票总仓位'
$ws.Range("F1").Value = '仓位占比'
$ws.Range("G1").Value = '持有市值(亿元)'
$ws.Range("H1").Value = '仓位排名'

# Copy the bordered/bold "header" style (cellXf index 2) from a donor
# cell onto every header cell in this sheet, B1:H1.
$wsQ4_2021.Range("B1").Copy()
$ws.Range("B1:H1").PasteSpecial(-4122)

# ---- data rows (rows 2..47), one per fund ----
# row 2: 161903
$ws.Cells.Item(2,2).Value = '161903'
$ws.Cells.Item(2,3).Value = '万家行业优选混合 (LOF)'
$ws.Cells.Item(2,4).NumberFormat = "@"
$ws.Cells.Item(2,4).Value = '112.51'
$ws.Cells.Item(2,4).Style = "Normal"
$ws.Cells.Item(2,5).NumberFormat = "@"
$ws.Cells.Item(2,5).Value = '91.94'
$ws.Cells.Item(2,5).Style = "Normal"
$ws.Cells.Item(2,6).NumberFormat = "@"
$ws.Cells.Item(2,6).Value = '5.99'
$ws.Cells.Item(2,6).Style = "Normal"
$ws.Cells.Item(2,7).NumberFormat = "@"
$ws.Cells.Item(2,7).Value = '6.7393'
$ws.Cells.Item(2,7).Style = "Normal"
$ws.Cells.Item(2,8).Value = 8

# row 3: 008120
$ws.Cells.Item(3,2).Value = '008120'
$ws.Cells.Item(3,3).Value = '万家自主创新混合A'
$ws.Cells.Item(3,4).NumberFormat = "@"
$ws.Cells.Item(3,4).Value = '31.86'
$ws.Cells.Item(3,4).Style = "Normal"
$ws.Cells.Item(3,5).NumberFormat = "@"
$ws.Cells.Item(3,5).Value = '93.78'
$ws.Cells.Item(3,5).Style = "Normal"
$ws.Cells.Item(3,6).NumberFormat = "@"
$ws.Cells.Item(3,6).Value = '7.19'
$ws.Cells.Item(3,6).Style = "Normal"
$ws.Cells.Item(3,7).NumberFormat = "@"
$ws.Cells.Item(3,7).Value = '2.2907'
$ws.Cells.Item(3,7).Style = "Normal"
$ws.Cells.Item(3,8).Value = 7

# row 4: 005311
$ws.Cells.Item(4,2).Value = '005311'
$ws.Cells.Item(4,3).Value = '万家经济新动能混合A'
$ws.Cells.Item(4,4).NumberFormat = "@"
$ws.Cells.Item(4,4).Value = '19.47'
$ws.Cells.Item(4,4).Style = "Normal"
$ws.Cells.Item(4,5).NumberFormat = "@"
$ws.Cells.Item(4,5).Value = '93.80'
$ws.Cells.Item(4,5).Style = "Normal"
$ws.Cells.Item(4,6).NumberFormat = "@"
$ws.Cells.Item(4,6).Value = '7.45'
$ws.Cells.Item(4,6).Style = "Normal"
$ws.Cells.Item(4,7).NumberFormat = "@"
$ws.Cells.Item(4,7).Value = '1.4505'
$ws.Cells.Item(4,7).Style = "Normal"
$ws.Cells.Item(4,8).Value = 8

# row 5: 159998
$ws.Cells.Item(5,2).Value = '159998'
$ws.Cells.Item(5,3).Value = '天弘中证计算机主题ETF'
$ws.Cells.Item(5,4).NumberFormat = "@"
$ws.Cells.Item(5,4).Value = '18.78'
$ws.Cells.Item(5,4).Style = "Normal"
$ws.Cells.Item(5,5).NumberFormat = "@"
$ws.Cells.Item(5,5).Value = '99.81'
$ws.Cells.Item(5,5).Style = "Normal"
$ws.Cells.Item(5,6).NumberFormat = "@"
$ws.Cells.Item(5,6).Value = '3.43'
$ws.Cells.Item(5,6).Style = "Normal"
$ws.Cells.Item(5,7).NumberFormat = "@"
$ws.Cells.Item(5,7).Value = '0.6442'
$ws.Cells.Item(5,7).Style = "Normal"
$ws.Cells.Item(5,8).Value = 5

# row 6: 159819
$ws.Cells.Item(6,2).Value = '159819'
$ws.Cells.Item(6,3).Value = '易方达中证人工智能主题ETF'
$ws.Cells.Item(6,4).NumberFormat = "@"
$ws.Cells.Item(6,4).Value = '14.51'
$ws.Cells.Item(6,4).Style = "Normal"
$ws.Cells.Item(6,5).NumberFormat = "@"
$ws.Cells.Item(6,5).Value = '99.17'
$ws.Cells.Item(6,5).Style = "Normal"
$ws.Cells.Item(6,6).NumberFormat = "@"
$ws.Cells.Item(6,6).Value = '3.80'
$ws.Cells.Item(6,6).Style = "Normal"
$ws.Cells.Item(6,7).NumberFormat = "@"
$ws.Cells.Item(6,7).Value = '0.5514'
$ws.Cells.Item(6,7).Style = "Normal"
$ws.Cells.Item(6,8).Value = 4

# row 7: 005312
$ws.Cells.Item(7,2).Value = '005312'
$ws.Cells.Item(7,3).Value = '万家经济新动能混合C'
$ws.Cells.Item(7,4).NumberFormat = "@"
$ws.Cells.Item(7,4).Value = '6.77'
$ws.Cells.Item(7,4).Style = "Normal"
$ws.Cells.Item(7,5).NumberFormat = "@"
$ws.Cells.Item(7,5).Value = '93.80'
$ws.Cells.Item(7,5).Style = "Normal"
$ws.Cells.Item(7,6).NumberFormat = "@"
$ws.Cells.Item(7,6).Value = '7.45'
$ws.Cells.Item(7,6).Style = "Normal"
$ws.Cells.Item(7,7).NumberFormat = "@"
$ws.Cells.Item(7,7).Value = '0.5044'
$ws.Cells.Item(7,7).Style = "Normal"
$ws.Cells.Item(7,8).Value = 8

# row 8: 000971
$ws.Cells.Item(8,2).Value = '000971'
$ws.Cells.Item(8,3).Value = '诺安新经济股票'
$ws.Cells.Item(8,4).NumberFormat = "@"
$ws.Cells.Item(8,4).Value = '15.43'
$ws.Cells.Item(8,4).Style = "Normal"
$ws.Cells.Item(8,5).NumberFormat = "@"
$ws.Cells.Item(8,5).Value = '82.95'
$ws.Cells.Item(8,5).Style = "Normal"
$ws.Cells.Item(8,6).NumberFormat = "@"
$ws.Cells.Item(8,6).Value = '3.07'
$ws.Cells.Item(8,6).Style = "Normal"
$ws.Cells.Item(8,7).NumberFormat = "@"
$ws.Cells.Item(8,7).Value = '0.4737'
$ws.Cells.Item(8,7).Style = "Normal"
$ws.Cells.Item(8,8).Value = 10

# row 9: 007639
$ws.Cells.Item(9,2).Value = '007639'
$ws.Cells.Item(9,3).Value = '汇添富3年封闭运作竞争优势灵活配置混合'
$ws.Cells.Item(9,4).NumberFormat = "@"
$ws.Cells.Item(9,4).Value = '13.07'
$ws.Cells.Item(9,4).Style = "Normal"
$ws.Cells.Item(9,5).NumberFormat = "@"
$ws.Cells.Item(9,5).Value = '60.50'
$ws.Cells.Item(9,5).Style = "Normal"
$ws.Cells.Item(9,6).NumberFormat = "@"
$ws.Cells.Item(9,6).Value = '2.45'
$ws.Cells.Item(9,6).Style = "Normal"
$ws.Cells.Item(9,7).NumberFormat = "@"
$ws.Cells.Item(9,7).Value = '0.3202'
$ws.Cells.Item(9,7).Style = "Normal"
$ws.Cells.Item(9,8).Value = 9

# row 10: 008633
$ws.Cells.Item(10,2).Value = '008633'
$ws.Cells.Item(10,3).Value = '万家科技创新混合A'
$ws.Cells.Item(10,4).NumberFormat = "@"
$ws.Cells.Item(10,4).Value = '3.75'
$ws.Cells.Item(10,4).Style = "Normal"
$ws.Cells.Item(10,5).NumberFormat = "@"
$ws.Cells.Item(10,5).Value = '93.27'
$ws.Cells.Item(10,5).Style = "Normal"
$ws.Cells.Item(10,6).NumberFormat = "@"
$ws.Cells.Item(10,6).Value = '7.53'
$ws.Cells.Item(10,6).Style = "Normal"
$ws.Cells.Item(10,7).NumberFormat = "@"
$ws.Cells.Item(10,7).Value = '0.2824'
$ws.Cells.Item(10,7).Style = "Normal"
$ws.Cells.Item(10,8).Value = 8

# row 11: 515070
$ws.Cells.Item(11,2).Value = '515070'
$ws.Cells.Item(11,3).Value = '华夏中证人工智能主题ETF'
$ws.Cells.Item(11,4).NumberFormat = "@"
$ws.Cells.Item(11,4).Value = '7.39'
$ws.Cells.Item(11,4).Style = "Normal"
$ws.Cells.Item(11,5).NumberFormat = "@"
$ws.Cells.Item(11,5).Value = '99.43'
$ws.Cells.Item(11,5).Style = "Normal"
$ws.Cells.Item(11,6).NumberFormat = "@"
$ws.Cells.Item(11,6).Value = '3.80'
$ws.Cells.Item(11,6).Style = "Normal"
$ws.Cells.Item(11,7).NumberFormat = "@"
$ws.Cells.Item(11,7).Value = '0.2808'
$ws.Cells.Item(11,7).Style = "Normal"
$ws.Cells.Item(11,8).Value = 4

# row 12: 512720
$ws.Cells.Item(12,2).Value = '512720'
$ws.Cells.Item(12,3).Value = '国泰中证计算机主题ETF'
$ws.Cells.Item(12,4).NumberFormat = "@"
$ws.Cells.Item(12,4).Value = '8.13'
$ws.Cells.Item(12,4).Style = "Normal"
$ws.Cells.Item(12,5).NumberFormat = "@"
$ws.Cells.Item(12,5).Value = '99.18'
$ws.Cells.Item(12,5).Style = "Normal"
$ws.Cells.Item(12,6).NumberFormat = "@"
$ws.Cells.Item(12,6).Value = '3.25'
$ws.Cells.Item(12,6).Style = "Normal"
$ws.Cells.Item(12,7).NumberFormat = "@"
$ws.Cells.Item(12,7).Value = '0.2642'
$ws.Cells.Item(12,7).Style = "Normal"
$ws.Cells.Item(12,8).Value = 6

# row 13: 161631
$ws.Cells.Item(13,2).Value = '161631'
$ws.Cells.Item(13,3).Value = '融通中证人工智能主题指数(LOF)A'
$ws.Cells.Item(13,4).NumberFormat = "@"
$ws.Cells.Item(13,4).Value = '5.66'
$ws.Cells.Item(13,4).Style = "Normal"
$ws.Cells.Item(13,5).NumberFormat = "@"
$ws.Cells.Item(13,5).Value = '94.47'
$ws.Cells.Item(13,5).Style = "Normal"
$ws.Cells.Item(13,6).NumberFormat = "@"
$ws.Cells.Item(13,6).Value = '3.66'
$ws.Cells.Item(13,6).Style = "Normal"
$ws.Cells.Item(13,7).NumberFormat = "@"
$ws.Cells.Item(13,7).Value = '0.2072'
$ws.Cells.Item(13,7).Style = "Normal"
$ws.Cells.Item(13,8).Value = 4

# row 14: 008121
$ws.Cells.Item(14,2).Value = '008121'
$ws.Cells.Item(14,3).Value = '万家自主创新混合C'
$ws.Cells.Item(14,4).NumberFormat = "@"
$ws.Cells.Item(14,4).Value = '2.80'
$ws.Cells.Item(14,4).Style = "Normal"
$ws.Cells.Item(14,5).NumberFormat = "@"
$ws.Cells.Item(14,5).Value = '93.78'
$ws.Cells.Item(14,5).Style = "Normal"
$ws.Cells.Item(14,6).NumberFormat = "@"
$ws.Cells.Item(14,6).Value = '7.19'
$ws.Cells.Item(14,6).Style = "Normal"
$ws.Cells.Item(14,7).NumberFormat = "@"
$ws.Cells.Item(14,7).Value = '0.2013'
$ws.Cells.Item(14,7).Style = "Normal"
$ws.Cells.Item(14,8).Value = 7

# row 15: 159852
$ws.Cells.Item(15,2).Value = '159852'
$ws.Cells.Item(15,3).Value = '嘉实中证软件服务ETF'
$ws.Cells.Item(15,4).NumberFormat = "@"
$ws.Cells.Item(15,4).Value = '2.78'
$ws.Cells.Item(15,4).Style = "Normal"
$ws.Cells.Item(15,5).NumberFormat = "@"
$ws.Cells.Item(15,5).Value = '98.71'
$ws.Cells.Item(15,5).Style = "Normal"
$ws.Cells.Item(15,6).NumberFormat = "@"
$ws.Cells.Item(15,6).Value = '5.87'
$ws.Cells.Item(15,6).Style = "Normal"
$ws.Cells.Item(15,7).NumberFormat = "@"
$ws.Cells.Item(15,7).Value = '0.1632'
$ws.Cells.Item(15,7).Style = "Normal"
$ws.Cells.Item(15,8).Value = 4

# row 16: 165523
$ws.Cells.Item(16,2).Value = '165523'
$ws.Cells.Item(16,3).Value = '信诚中证信息安全指数（LOF）'
$ws.Cells.Item(16,4).NumberFormat = "@"
$ws.Cells.Item(16,4).Value = '4.04'
$ws.Cells.Item(16,4).Style = "Normal"
$ws.Cells.Item(16,5).NumberFormat = "@"
$ws.Cells.Item(16,5).Value = '93.97'
$ws.Cells.Item(16,5).Style = "Normal"
$ws.Cells.Item(16,6).NumberFormat = "@"
$ws.Cells.Item(16,6).Value = '3.98'
$ws.Cells.Item(16,6).Style = "Normal"
$ws.Cells.Item(16,7).NumberFormat = "@"
$ws.Cells.Item(16,7).Value = '0.1608'
$ws.Cells.Item(16,7).Style = "Normal"
$ws.Cells.Item(16,8).Value = 4

# row 17: 360011
$ws.Cells.Item(17,2).Value = '360011'
$ws.Cells.Item(17,3).Value = '光大保德信动态优选混合'
$ws.Cells.Item(17,4).NumberFormat = "@"
$ws.Cells.Item(17,4).Value = '2.56'
$ws.Cells.Item(17,4).Style = "Normal"
$ws.Cells.Item(17,5).NumberFormat = "@"
$ws.Cells.Item(17,5).Value = '71.83'
$ws.Cells.Item(17,5).Style = "Normal"
$ws.Cells.Item(17,6).NumberFormat = "@"
$ws.Cells.Item(17,6).Value = '5.90'
$ws.Cells.Item(17,6).Style = "Normal"
$ws.Cells.Item(17,7).NumberFormat = "@"
$ws.Cells.Item(17,7).Value = '0.1510'
$ws.Cells.Item(17,7).Style = "Normal"
$ws.Cells.Item(17,8).Value = 6

# row 18: 161628
$ws.Cells.Item(18,2).Value = '161628'
$ws.Cells.Item(18,3).Value = '融通中证云计算与大数据主题指数（LOF）'
$ws.Cells.Item(18,4).NumberFormat = "@"
$ws.Cells.Item(18,4).Value = '2.71'
$ws.Cells.Item(18,4).Style = "Normal"
$ws.Cells.Item(18,5).NumberFormat = "@"
$ws.Cells.Item(18,5).Value = '92.74'
$ws.Cells.Item(18,5).Style = "Normal"
$ws.Cells.Item(18,6).NumberFormat = "@"
$ws.Cells.Item(18,6).Value = '5.43'
$ws.Cells.Item(18,6).Style = "Normal"
$ws.Cells.Item(18,7).NumberFormat = "@"
$ws.Cells.Item(18,7).Value = '0.1472'
$ws.Cells.Item(18,7).Style = "Normal"
$ws.Cells.Item(18,8).Value = 3

# row 19: 515230
$ws.Cells.Item(19,2).Value = '515230'
$ws.Cells.Item(19,3).Value = '国泰中证全指软件ETF'
$ws.Cells.Item(19,4).NumberFormat = "@"
$ws.Cells.Item(19,4).Value = '2.53'
$ws.Cells.Item(19,4).Style = "Normal"
$ws.Cells.Item(19,5).NumberFormat = "@"
$ws.Cells.Item(19,5).Value = '98.89'
$ws.Cells.Item(19,5).Style = "Normal"
$ws.Cells.Item(19,6).NumberFormat = "@"
$ws.Cells.Item(19,6).Value = '5.44'
$ws.Cells.Item(19,6).Style = "Normal"
$ws.Cells.Item(19,7).NumberFormat = "@"
$ws.Cells.Item(19,7).Value = '0.1376'
$ws.Cells.Item(19,7).Style = "Normal"
$ws.Cells.Item(19,8).Value = 4

# row 20: 515400
$ws.Cells.Item(20,2).Value = '515400'
$ws.Cells.Item(20,3).Value = '富国中证大数据产业ETF'
$ws.Cells.Item(20,4).NumberFormat = "@"
$ws.Cells.Item(20,4).Value = '2.29'
$ws.Cells.Item(20,4).Style = "Normal"
$ws.Cells.Item(20,5).NumberFormat = "@"
$ws.Cells.Item(20,5).Value = '99.33'
$ws.Cells.Item(20,5).Style = "Normal"
$ws.Cells.Item(20,6).NumberFormat = "@"
$ws.Cells.Item(20,6).Value = '6.00'
$ws.Cells.Item(20,6).Style = "Normal"
$ws.Cells.Item(20,7).NumberFormat = "@"
$ws.Cells.Item(20,7).Value = '0.1374'
$ws.Cells.Item(20,7).Style = "Normal"
$ws.Cells.Item(20,8).Value = 4

# row 21: 515980
$ws.Cells.Item(21,2).Value = '515980'
$ws.Cells.Item(21,3).Value = '华富中证人工智能产业ETF'
$ws.Cells.Item(21,4).NumberFormat = "@"
$ws.Cells.Item(21,4).Value = '2.93'
$ws.Cells.Item(21,4).Style = "Normal"
$ws.Cells.Item(21,5).NumberFormat = "@"
$ws.Cells.Item(21,5).Value = '99.14'
$ws.Cells.Item(21,5).Style = "Normal"
$ws.Cells.Item(21,6).NumberFormat = "@"
$ws.Cells.Item(21,6).Value = '4.14'
$ws.Cells.Item(21,6).Style = "Normal"
$ws.Cells.Item(21,7).NumberFormat = "@"
$ws.Cells.Item(21,7).Value = '0.1213'
$ws.Cells.Item(21,7).Style = "Normal"
$ws.Cells.Item(21,8).Value = 4

# row 22: 516510
$ws.Cells.Item(22,2).Value = '516510'
$ws.Cells.Item(22,3).Value = '易方达中证云计算与大数据主题ETF'
$ws.Cells.Item(22,4).NumberFormat = "@"
$ws.Cells.Item(22,4).Value = '2.10'
$ws.Cells.Item(22,4).Style = "Normal"
$ws.Cells.Item(22,5).NumberFormat = "@"
$ws.Cells.Item(22,5).Value = '98.62'
$ws.Cells.Item(22,5).Style = "Normal"
$ws.Cells.Item(22,6).NumberFormat = "@"
$ws.Cells.Item(22,6).Value = '5.70'
$ws.Cells.Item(22,6).Style = "Normal"
$ws.Cells.Item(22,7).NumberFormat = "@"
$ws.Cells.Item(22,7).Value = '0.1197'
$ws.Cells.Item(22,7).Style = "Normal"
$ws.Cells.Item(22,8).Value = 3

# row 23: 560660
$ws.Cells.Item(23,2).Value = '560660'
$ws.Cells.Item(23,3).Value = '新华中证云计算50交易型开放式指数证券投资基金'
$ws.Cells.Item(23,4).NumberFormat = "@"
$ws.Cells.Item(23,4).Value = '1.75'
$ws.Cells.Item(23,4).Style = "Normal"
$ws.Cells.Item(23,5).NumberFormat = "@"
$ws.Cells.Item(23,5).Value = '97.21'
$ws.Cells.Item(23,5).Style = "Normal"
$ws.Cells.Item(23,6).NumberFormat = "@"
$ws.Cells.Item(23,6).Value = '6.80'
$ws.Cells.Item(23,6).Style = "Normal"
$ws.Cells.Item(23,7).NumberFormat = "@"
$ws.Cells.Item(23,7).Value = '0.1190'
$ws.Cells.Item(23,7).Style = "Normal"
$ws.Cells.Item(23,8).Value = 3

# row 24: 007854
$ws.Cells.Item(24,2).Value = '007854'
$ws.Cells.Item(24,3).Value = '光大保德信景气先锋混合'
$ws.Cells.Item(24,4).NumberFormat = "@"
$ws.Cells.Item(24,4).Value = '1.82'
$ws.Cells.Item(24,4).Style = "Normal"
$ws.Cells.Item(24,5).NumberFormat = "@"
$ws.Cells.Item(24,5).Value = '72.20'
$ws.Cells.Item(24,5).Style = "Normal"
$ws.Cells.Item(24,6).NumberFormat = "@"
$ws.Cells.Item(24,6).Value = '5.98'
$ws.Cells.Item(24,6).Style = "Normal"
$ws.Cells.Item(24,7).NumberFormat = "@"
$ws.Cells.Item(24,7).Value = '0.1088'
$ws.Cells.Item(24,7).Style = "Normal"
$ws.Cells.Item(24,8).Value = 5

# row 25: 008634
$ws.Cells.Item(25,2).Value = '008634'
$ws.Cells.Item(25,3).Value = '万家科技创新混合C'
$ws.Cells.Item(25,4).NumberFormat = "@"
$ws.Cells.Item(25,4).Value = '1.18'
$ws.Cells.Item(25,4).Style = "Normal"
$ws.Cells.Item(25,5).NumberFormat = "@"
$ws.Cells.Item(25,5).Value = '93.27'
$ws.Cells.Item(25,5).Style = "Normal"
$ws.Cells.Item(25,6).NumberFormat = "@"
$ws.Cells.Item(25,6).Value = '7.53'
$ws.Cells.Item(25,6).Style = "Normal"
$ws.Cells.Item(25,7).NumberFormat = "@"
$ws.Cells.Item(25,7).Value = '0.0889'
$ws.Cells.Item(25,7).Style = "Normal"
$ws.Cells.Item(25,8).Value = 8

# row 26: 512930
$ws.Cells.Item(26,2).Value = '512930'
$ws.Cells.Item(26,3).Value = '平安中证人工智能主题ETF'
$ws.Cells.Item(26,4).NumberFormat = "@"
$ws.Cells.Item(26,4).Value = '2.31'
$ws.Cells.Item(26,4).Style = "Normal"
$ws.Cells.Item(26,5).NumberFormat = "@"
$ws.Cells.Item(26,5).Value = '98.49'
$ws.Cells.Item(26,5).Style = "Normal"
$ws.Cells.Item(26,6).NumberFormat = "@"
$ws.Cells.Item(26,6).Value = '3.73'
$ws.Cells.Item(26,6).Style = "Normal"
$ws.Cells.Item(26,7).NumberFormat = "@"
$ws.Cells.Item(26,7).Value = '0.0862'
$ws.Cells.Item(26,7).Style = "Normal"
$ws.Cells.Item(26,8).Value = 4

# row 27: 159613
$ws.Cells.Item(27,2).Value = '159613'
$ws.Cells.Item(27,3).Value = '嘉实中证信息安全主题ETF'
$ws.Cells.Item(27,4).NumberFormat = "@"
$ws.Cells.Item(27,4).Value = '1.82'
$ws.Cells.Item(27,4).Style = "Normal"
$ws.Cells.Item(27,5).NumberFormat = "@"
$ws.Cells.Item(27,5).Value = '99.01'
$ws.Cells.Item(27,5).Style = "Normal"
$ws.Cells.Item(27,6).NumberFormat = "@"
$ws.Cells.Item(27,6).Value = '4.18'
$ws.Cells.Item(27,6).Style = "Normal"
$ws.Cells.Item(27,7).NumberFormat = "@"
$ws.Cells.Item(27,7).Value = '0.0761'
$ws.Cells.Item(27,7).Style = "Normal"
$ws.Cells.Item(27,8).Value = 4

# row 28: 159899
$ws.Cells.Item(28,2).Value = '159899'
$ws.Cells.Item(28,3).Value = '招商中证全指软件交易型开放式指数证券投资基金'
$ws.Cells.Item(28,4).NumberFormat = "@"
$ws.Cells.Item(28,4).Value = '1.26'
$ws.Cells.Item(28,4).Style = "Normal"
$ws.Cells.Item(28,5).NumberFormat = "@"
$ws.Cells.Item(28,5).Value = '98.64'
$ws.Cells.Item(28,5).Style = "Normal"
$ws.Cells.Item(28,6).NumberFormat = "@"
$ws.Cells.Item(28,6).Value = '5.54'
$ws.Cells.Item(28,6).Style = "Normal"
$ws.Cells.Item(28,7).NumberFormat = "@"
$ws.Cells.Item(28,7).Value = '0.0698'
$ws.Cells.Item(28,7).Style = "Normal"
$ws.Cells.Item(28,8).Value = 4

# row 29: 501063
$ws.Cells.Item(29,2).Value = '501063'
$ws.Cells.Item(29,3).Value = '汇添富悦享定期开放混合'
$ws.Cells.Item(29,4).NumberFormat = "@"
$ws.Cells.Item(29,4).Value = '2.18'
$ws.Cells.Item(29,4).Style = "Normal"
$ws.Cells.Item(29,5).NumberFormat = "@"
$ws.Cells.Item(29,5).Value = '60.21'
$ws.Cells.Item(29,5).Style = "Normal"
$ws.Cells.Item(29,6).NumberFormat = "@"
$ws.Cells.Item(29,6).Value = '2.63'
$ws.Cells.Item(29,6).Style = "Normal"
$ws.Cells.Item(29,7).NumberFormat = "@"
$ws.Cells.Item(29,7).Value = '0.0573'
$ws.Cells.Item(29,7).Style = "Normal"
$ws.Cells.Item(29,8).Value = 10

# row 30: 000589
$ws.Cells.Item(30,2).Value = '000589'
$ws.Cells.Item(30,3).Value = '光大保德信银发商机主题混合'
$ws.Cells.Item(30,4).NumberFormat = "@"
$ws.Cells.Item(30,4).Value = '1.30'
$ws.Cells.Item(30,4).Style = "Normal"
$ws.Cells.Item(30,5).NumberFormat = "@"
$ws.Cells.Item(30,5).Value = '87.65'
$ws.Cells.Item(30,5).Style = "Normal"
$ws.Cells.Item(30,6).NumberFormat = "@"
$ws.Cells.Item(30,6).Value = '3.44'
$ws.Cells.Item(30,6).Style = "Normal"
$ws.Cells.Item(30,7).NumberFormat = "@"
$ws.Cells.Item(30,7).Value = '0.0447'
$ws.Cells.Item(30,7).Style = "Normal"
$ws.Cells.Item(30,8).Value = 7

# row 31: 516630
$ws.Cells.Item(31,2).Value = '516630'
$ws.Cells.Item(31,3).Value = '华夏中证云计算与大数据主题ETF'
$ws.Cells.Item(31,4).NumberFormat = "@"
$ws.Cells.Item(31,4).Value = '0.76'
$ws.Cells.Item(31,4).Style = "Normal"
$ws.Cells.Item(31,5).NumberFormat = "@"
$ws.Cells.Item(31,5).Value = '98.87'
$ws.Cells.Item(31,5).Style = "Normal"
$ws.Cells.Item(31,6).NumberFormat = "@"
$ws.Cells.Item(31,6).Value = '5.71'
$ws.Cells.Item(31,6).Style = "Normal"
$ws.Cells.Item(31,7).NumberFormat = "@"
$ws.Cells.Item(31,7).Value = '0.0434'
$ws.Cells.Item(31,7).Style = "Normal"
$ws.Cells.Item(31,8).Value = 3

# row 32: 517200
$ws.Cells.Item(32,2).Value = '517200'
$ws.Cells.Item(32,3).Value = '嘉实中证沪港深互联网ETF'
$ws.Cells.Item(32,4).NumberFormat = "@"
$ws.Cells.Item(32,4).Value = '1.71'
$ws.Cells.Item(32,4).Style = "Normal"
$ws.Cells.Item(32,5).NumberFormat = "@"
$ws.Cells.Item(32,5).Value = '98.76'
$ws.Cells.Item(32,5).Style = "Normal"
$ws.Cells.Item(32,6).NumberFormat = "@"
$ws.Cells.Item(32,6).Value = '2.46'
$ws.Cells.Item(32,6).Style = "Normal"
$ws.Cells.Item(32,7).NumberFormat = "@"
$ws.Cells.Item(32,7).Value = '0.0421'
$ws.Cells.Item(32,7).Style = "Normal"
$ws.Cells.Item(32,8).Value = 10

# row 33: 168701
$ws.Cells.Item(33,2).Value = '168701'
$ws.Cells.Item(33,3).Value = '合煦智远国证香蜜湖金融科技指数(LOF)A'
$ws.Cells.Item(33,4).NumberFormat = "@"
$ws.Cells.Item(33,4).Value = '0.90'
$ws.Cells.Item(33,4).Style = "Normal"
$ws.Cells.Item(33,5).NumberFormat = "@"
$ws.Cells.Item(33,5).Value = '93.15'
$ws.Cells.Item(33,5).Style = "Normal"
$ws.Cells.Item(33,6).NumberFormat = "@"
$ws.Cells.Item(33,6).Value = '4.41'
$ws.Cells.Item(33,6).Style = "Normal"
$ws.Cells.Item(33,7).NumberFormat = "@"
$ws.Cells.Item(33,7).Value = '0.0397'
$ws.Cells.Item(33,7).Style = "Normal"
$ws.Cells.Item(33,8).Value = 5

# row 34: 516000
$ws.Cells.Item(34,2).Value = '516000'
$ws.Cells.Item(34,3).Value = '华夏中证大数据产业ETF'
$ws.Cells.Item(34,4).NumberFormat = "@"
$ws.Cells.Item(34,4).Value = '0.58'
$ws.Cells.Item(34,4).Style = "Normal"
$ws.Cells.Item(34,5).NumberFormat = "@"
$ws.Cells.Item(34,5).Value = '97.61'
$ws.Cells.Item(34,5).Style = "Normal"
$ws.Cells.Item(34,6).NumberFormat = "@"
$ws.Cells.Item(34,6).Value = '5.92'
$ws.Cells.Item(34,6).Style = "Normal"
$ws.Cells.Item(34,7).NumberFormat = "@"
$ws.Cells.Item(34,7).Value = '0.0343'
$ws.Cells.Item(34,7).Style = "Normal"
$ws.Cells.Item(34,8).Value = 4

# row 35: 012371
$ws.Cells.Item(35,2).Value = '012371'
$ws.Cells.Item(35,3).Value = '西藏东财中证沪港深互联网指数型发起式证券投资基金A'
$ws.Cells.Item(35,4).NumberFormat = "@"
$ws.Cells.Item(35,4).Value = '1.11'
$ws.Cells.Item(35,4).Style = "Normal"
$ws.Cells.Item(35,5).NumberFormat = "@"
$ws.Cells.Item(35,5).Value = '95.04'
$ws.Cells.Item(35,5).Style = "Normal"
$ws.Cells.Item(35,6).NumberFormat = "@"
$ws.Cells.Item(35,6).Value = '2.37'
$ws.Cells.Item(35,6).Style = "Normal"
$ws.Cells.Item(35,7).NumberFormat = "@"
$ws.Cells.Item(35,7).Value = '0.0263'
$ws.Cells.Item(35,7).Style = "Normal"
$ws.Cells.Item(35,8).Value = 10

# row 36: 011839
$ws.Cells.Item(36,2).Value = '011839'
$ws.Cells.Item(36,3).Value = '天弘中证人工智能指数A'
$ws.Cells.Item(36,4).NumberFormat = "@"
$ws.Cells.Item(36,4).Value = '0.72'
$ws.Cells.Item(36,4).Style = "Normal"
$ws.Cells.Item(36,5).NumberFormat = "@"
$ws.Cells.Item(36,5).Value = '94.99'
$ws.Cells.Item(36,5).Style = "Normal"
$ws.Cells.Item(36,6).NumberFormat = "@"
$ws.Cells.Item(36,6).Value = '3.65'
$ws.Cells.Item(36,6).Style = "Normal"
$ws.Cells.Item(36,7).NumberFormat = "@"
$ws.Cells.Item(36,7).Value = '0.0263'
$ws.Cells.Item(36,7).Style = "Normal"
$ws.Cells.Item(36,8).Value = 4

# row 37: 159890
$ws.Cells.Item(37,2).Value = '159890'
$ws.Cells.Item(37,3).Value = '招商中证云计算与大数据主题ETF'
$ws.Cells.Item(37,4).NumberFormat = "@"
$ws.Cells.Item(37,4).Value = '0.47'
$ws.Cells.Item(37,4).Style = "Normal"
$ws.Cells.Item(37,5).NumberFormat = "@"
$ws.Cells.Item(37,5).Value = '97.81'
$ws.Cells.Item(37,5).Style = "Normal"
$ws.Cells.Item(37,6).NumberFormat = "@"
$ws.Cells.Item(37,6).Value = '5.49'
$ws.Cells.Item(37,6).Style = "Normal"
$ws.Cells.Item(37,7).NumberFormat = "@"
$ws.Cells.Item(37,7).Value = '0.0258'
$ws.Cells.Item(37,7).Style = "Normal"
$ws.Cells.Item(37,8).Value = 3

# row 38: 011840
$ws.Cells.Item(38,2).Value = '011840'
$ws.Cells.Item(38,3).Value = '天弘中证人工智能指数C'
$ws.Cells.Item(38,4).NumberFormat = "@"
$ws.Cells.Item(38,4).Value = '0.58'
$ws.Cells.Item(38,4).Style = "Normal"
$ws.Cells.Item(38,5).NumberFormat = "@"
$ws.Cells.Item(38,5).Value = '94.99'
$ws.Cells.Item(38,5).Style = "Normal"
$ws.Cells.Item(38,6).NumberFormat = "@"
$ws.Cells.Item(38,6).Value = '3.65'
$ws.Cells.Item(38,6).Style = "Normal"
$ws.Cells.Item(38,7).NumberFormat = "@"
$ws.Cells.Item(38,7).Value = '0.0212'
$ws.Cells.Item(38,7).Style = "Normal"
$ws.Cells.Item(38,8).Value = 4

# row 39: 009239
$ws.Cells.Item(39,2).Value = '009239'
$ws.Cells.Item(39,3).Value = '融通中证人工智能主题指数(LOF)C'
$ws.Cells.Item(39,4).NumberFormat = "@"
$ws.Cells.Item(39,4).Value = '0.51'
$ws.Cells.Item(39,4).Style = "Normal"
$ws.Cells.Item(39,5).NumberFormat = "@"
$ws.Cells.Item(39,5).Value = '94.47'
$ws.Cells.Item(39,5).Style = "Normal"
$ws.Cells.Item(39,6).NumberFormat = "@"
$ws.Cells.Item(39,6).Value = '3.66'
$ws.Cells.Item(39,6).Style = "Normal"
$ws.Cells.Item(39,7).NumberFormat = "@"
$ws.Cells.Item(39,7).Value = '0.0187'
$ws.Cells.Item(39,7).Style = "Normal"
$ws.Cells.Item(39,8).Value = 4

# row 40: 516700
$ws.Cells.Item(40,2).Value = '516700'
$ws.Cells.Item(40,3).Value = '华宝中证大数据产业交易型开放式指数证券投资基金'
$ws.Cells.Item(40,4).NumberFormat = "@"
$ws.Cells.Item(40,4).Value = '0.31'
$ws.Cells.Item(40,4).Style = "Normal"
$ws.Cells.Item(40,5).NumberFormat = "@"
$ws.Cells.Item(40,5).Value = '96.93'
$ws.Cells.Item(40,5).Style = "Normal"
$ws.Cells.Item(40,6).NumberFormat = "@"
$ws.Cells.Item(40,6).Value = '5.88'
$ws.Cells.Item(40,6).Style = "Normal"
$ws.Cells.Item(40,7).NumberFormat = "@"
$ws.Cells.Item(40,7).Value = '0.0182'
$ws.Cells.Item(40,7).Style = "Normal"
$ws.Cells.Item(40,8).Value = 4

# row 41: 001780
$ws.Cells.Item(41,2).Value = '001780'
$ws.Cells.Item(41,3).Value = '诺安改革趋势灵活配置混合'
$ws.Cells.Item(41,4).NumberFormat = "@"
$ws.Cells.Item(41,4).Value = '0.46'
$ws.Cells.Item(41,4).Style = "Normal"
$ws.Cells.Item(41,5).NumberFormat = "@"
$ws.Cells.Item(41,5).Value = '68.34'
$ws.Cells.Item(41,5).Style = "Normal"
$ws.Cells.Item(41,6).NumberFormat = "@"
$ws.Cells.Item(41,6).Value = '3.72'
$ws.Cells.Item(41,6).Style = "Normal"
$ws.Cells.Item(41,7).NumberFormat = "@"
$ws.Cells.Item(41,7).Value = '0.0171'
$ws.Cells.Item(41,7).Style = "Normal"
$ws.Cells.Item(41,8).Value = 9

# row 42: 517800
$ws.Cells.Item(42,2).Value = '517800'
$ws.Cells.Item(42,3).Value = '方正富邦中证沪港深人工智能50交易型开放式指数证券投资基金'
$ws.Cells.Item(42,4).NumberFormat = "@"
$ws.Cells.Item(42,4).Value = '0.44'
$ws.Cells.Item(42,4).Style = "Normal"
$ws.Cells.Item(42,5).NumberFormat = "@"
$ws.Cells.Item(42,5).Value = '90.64'
$ws.Cells.Item(42,5).Style = "Normal"
$ws.Cells.Item(42,6).NumberFormat = "@"
$ws.Cells.Item(42,6).Value = '3.12'
$ws.Cells.Item(42,6).Style = "Normal"
$ws.Cells.Item(42,7).NumberFormat = "@"
$ws.Cells.Item(42,7).Value = '0.0137'
$ws.Cells.Item(42,7).Style = "Normal"
$ws.Cells.Item(42,8).Value = 9

# row 43: 012372
$ws.Cells.Item(43,2).Value = '012372'
$ws.Cells.Item(43,3).Value = '西藏东财中证沪港深互联网指数型发起式证券投资基金C'
$ws.Cells.Item(43,4).NumberFormat = "@"
$ws.Cells.Item(43,4).Value = '0.51'
$ws.Cells.Item(43,4).Style = "Normal"
$ws.Cells.Item(43,5).NumberFormat = "@"
$ws.Cells.Item(43,5).Value = '95.04'
$ws.Cells.Item(43,5).Style = "Normal"
$ws.Cells.Item(43,6).NumberFormat = "@"
$ws.Cells.Item(43,6).Value = '2.37'
$ws.Cells.Item(43,6).Style = "Normal"
$ws.Cells.Item(43,7).NumberFormat = "@"
$ws.Cells.Item(43,7).Value = '0.0121'
$ws.Cells.Item(43,7).Style = "Normal"
$ws.Cells.Item(43,8).Value = 10

# row 44: 168702
$ws.Cells.Item(44,2).Value = '168702'
$ws.Cells.Item(44,3).Value = '合煦智远国证香蜜湖金融科技指数(LOF)C'
$ws.Cells.Item(44,4).NumberFormat = "@"
$ws.Cells.Item(44,4).Value = '0.22'
$ws.Cells.Item(44,4).Style = "Normal"
$ws.Cells.Item(44,5).NumberFormat = "@"
$ws.Cells.Item(44,5).Value = '93.15'
$ws.Cells.Item(44,5).Style = "Normal"
$ws.Cells.Item(44,6).NumberFormat = "@"
$ws.Cells.Item(44,6).Value = '4.41'
$ws.Cells.Item(44,6).Style = "Normal"
$ws.Cells.Item(44,7).NumberFormat = "@"
$ws.Cells.Item(44,7).Value = '0.0097'
$ws.Cells.Item(44,7).Style = "Normal"
$ws.Cells.Item(44,8).Value = 5

# row 45: 014543
$ws.Cells.Item(45,2).Value = '014543'
$ws.Cells.Item(45,3).Value = '汇添富中证沪港深云计算产业指数A'
$ws.Cells.Item(45,4).NumberFormat = "@"
$ws.Cells.Item(45,4).Value = '0.15'
$ws.Cells.Item(45,4).Style = "Normal"
$ws.Cells.Item(45,5).NumberFormat = "@"
$ws.Cells.Item(45,5).Value = '92.49'
$ws.Cells.Item(45,5).Style = "Normal"
$ws.Cells.Item(45,6).NumberFormat = "@"
$ws.Cells.Item(45,6).Value = '4.62'
$ws.Cells.Item(45,6).Style = "Normal"
$ws.Cells.Item(45,7).NumberFormat = "@"
$ws.Cells.Item(45,7).Value = '0.0069'
$ws.Cells.Item(45,7).Style = "Normal"
$ws.Cells.Item(45,8).Value = 4

# row 46: 159702
$ws.Cells.Item(46,2).Value = '159702'
$ws.Cells.Item(46,3).Value = '汇添富中证人工智能主题交易型开放式指数证券投资基金'
$ws.Cells.Item(46,4).NumberFormat = "@"
$ws.Cells.Item(46,4).Value = '0.11'
$ws.Cells.Item(46,4).Style = "Normal"
$ws.Cells.Item(46,5).NumberFormat = "@"
$ws.Cells.Item(46,5).Value = '99.85'
$ws.Cells.Item(46,5).Style = "Normal"
$ws.Cells.Item(46,6).NumberFormat = "@"
$ws.Cells.Item(46,6).Value = '3.82'
$ws.Cells.Item(46,6).Style = "Normal"
$ws.Cells.Item(46,7).NumberFormat = "@"
$ws.Cells.Item(46,7).Value = '0.0042'
$ws.Cells.Item(46,7).Style = "Normal"
$ws.Cells.Item(46,8).Value = 4

# row 47: 014544
$ws.Cells.Item(47,2).Value = '014544'
$ws.Cells.Item(47,3).Value = '汇添富中证沪港深云计算产业指数C'
$ws.Cells.Item(47,4).NumberFormat = "@"
$ws.Cells.Item(47,4).Value = '0.04'
$ws.Cells.Item(47,4).Style = "Normal"
$ws.Cells.Item(47,5).NumberFormat = "@"
$ws.Cells.Item(47,5).Value = '92.49'
$ws.Cells.Item(47,5).Style = "Normal"
$ws.Cells.Item(47,6).NumberFormat = "@"
$ws.Cells.Item(47,6).Value = '4.62'
$ws.Cells.Item(47,6).Style = "Normal"
$ws.Cells.Item(47,7).NumberFormat = "@"
$ws.Cells.Item(47,7).Value = '0.0018'
$ws.Cells.Item(47,7).Style = "Normal"
$ws.Cells.Item(47,8).Value = 4

# column A: 0-indexed sequence number, style copied from the donor's
# column-A cell (same cellXf index 2 as the header).
$wsQ4_2021.Range("A2").Copy()
$ws.Range("A2:A47").PasteSpecial(-4122)
$ws.Cells.Item(2,1).Value = 0
$ws.Cells.Item(3,1).Value = 1
$ws.Cells.Item(4,1).Value = 2
$ws.Cells.Item(5,1).Value = 3
$ws.Cells.Item(6,1).Value = 4
$ws.Cells.Item(7,1).Value = 5
$ws.Cells.Item(8,1).Value = 6
$ws.Cells.Item(9,1).Value = 7
$ws.Cells.Item(10,1).Value = 8
$ws.Cells.Item(11,1).Value = 9
$ws.Cells.Item(12,1).Value = 10
$ws.Cells.Item(13,1).Value = 11
$ws.Cells.Item(14,1).Value = 12
$ws.Cells.Item(15,1).Value = 13
$ws.Cells.Item(16,1).Value = 14
$ws.Cells.Item(17,1).Value = 15
$ws.Cells.Item(18,1).Value = 16
$ws.Cells.Item(19,1).Value = 17
$ws.Cells.Item(20,1).Value = 18
$ws.Cells.Item(21,1).Value = 19
$ws.Cells.Item(22,1).Value = 20
$ws.Cells.Item(23,1).Value = 21
$ws.Cells.Item(24,1).Value = 22
$ws.Cells.Item(25,1).Value = 23
$ws.Cells.Item(26,1).Value = 24
$ws.Cells.Item(27,1).Value = 25
$ws.Cells.Item(28,1).Value = 26
$ws.Cells.Item(29,1).Value = 27
$ws.Cells.Item(30,1).Value = 28
$ws.Cells.Item(31,1).Value = 29
$ws.Cells.Item(32,1).Value = 30
$ws.Cells.Item(33,1).Value = 31
$ws.Cells.Item(34,1).Value = 32
$ws.Cells.Item(35,1).Value = 33
$ws.Cells.Item(36,1).Value = 34
$ws.Cells.Item(37,1).Value = 35
$ws.Cells.Item(38,1).Value = 36
$ws.Cells.Item(39,1).Value = 37
$ws.Cells.Item(40,1).Value = 38
$ws.Cells.Item(41,1).Value = 39
$ws.Cells.Item(42,1).Value = 40
$ws.Cells.Item(43,1).Value = 41
$ws.Cells.Item(44,1).Value = 42
$ws.Cells.Item(45,1).Value = 43
$ws.Cells.Item(46,1).Value = 44
$ws.Cells.Item(47,1).Value = 45

# ---------------------------------------------------------------------
# Step 2: append a brand-new worksheet after "2022-Q1"; this becomes the
# new "总计" summary sheet (sheetId 7).
# ---------------------------------------------------------------------
$wsNewTotal = $wb.Worksheets.Add($null, $ws)
$wsNewTotal.Name = "总计"

$wsNewTotal.Range("B1").Value = '日期'
$wsNewTotal.Range("C1").Value = '持有数量(只)'
$wsNewTotal.Range("D1").Value = '持有市值(亿元)'

# Re-use the same header style for the new summary sheet too.
$wsQ4_2021.Range("B1").Copy()
$wsNewTotal.Range("B1:D1").PasteSpecial(-4122)

# ---- data rows (rows 2..7): 2022-Q1 followed by the five existing quarters ----
# row 2: 2022-Q1
$wsNewTotal.Cells.Item(2,2).Value = '2022-Q1'
$wsNewTotal.Cells.Item(2,3).Value = 46
$wsNewTotal.Cells.Item(2,4).Value = 16.36

# row 3: 2021-Q4
$wsNewTotal.Cells.Item(3,2).Value = '2021-Q4'
$wsNewTotal.Cells.Item(3,3).Value = 82
$wsNewTotal.Cells.Item(3,4).Value = 32.54

# row 4: 2021-Q3
$wsNewTotal.Cells.Item(4,2).Value = '2021-Q3'
$wsNewTotal.Cells.Item(4,3).Value = 75
$wsNewTotal.Cells.Item(4,4).Value = 38.61

# row 5: 2021-Q2
$wsNewTotal.Cells.Item(5,2).Value = '2021-Q2'
$wsNewTotal.Cells.Item(5,3).Value = 55
$wsNewTotal.Cells.Item(5,4).Value = 35.87

# row 6: 2021-Q1
$wsNewTotal.Cells.Item(6,2).Value = '2021-Q1'
$wsNewTotal.Cells.Item(6,3).Value = 78
$wsNewTotal.Cells.Item(6,4).Value = 54

# row 7: 2020-Q4
$wsNewTotal.Cells.Item(7,2).Value = '2020-Q4'
$wsNewTotal.Cells.Item(7,3).Value = 96
$wsNewTotal.Cells.Item(7,4).Value = 83.36

# column A: 0-indexed sequence number, same donor style as above.
$wsQ4_2021.Range("A2").Copy()
$wsNewTotal.Range("A2:A7").PasteSpecial(-4122)
$wsNewTotal.Cells.Item(2,1).Value = 0
$wsNewTotal.Cells.Item(3,1).Value = 1
$wsNewTotal.Cells.Item(4,1).Value = 2
$wsNewTotal.Cells.Item(5,1).Value = 3
$wsNewTotal.Cells.Item(6,1).Value = 4
$wsNewTotal.Cells.Item(7,1).Value = 5
